$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Kerdesek")
$ws2 = $wb.Worksheets.Item("Valaszok")

# --- Fill sheet1 (Kerdesek) column A rows 1-24 with questions ---
$ws1.Range("A1").Value = "Kérdés"
$ws1.Range("A2").Value = "Két szomszéd lakik egymás mellett, mindenkit látnak, csak egymást nem láthatják. Mi az?"
$ws1.Range("A3").Value = "Veled megyen, nincs teste, napsütésben fekete."
$ws1.Range("A4").Value = "Fent lakom az égen, melegít a fényem, sugárból van bajszom, este van, ha alszom."
$ws1.Range("A5").Value = "Édesanyád gyermeke, de neked nem testvéred, ki az?"
$ws1.Range("A6").Value = "Szereted vagy nem szereted, ha megeszed, megkönnyezed."
$ws1.Range("A7").Value = "Két kezemmel füled fogom, s ott csücsülök az orrodon."
$ws1.Range("A8").Value = "Csak a tiéd, mégis mások használják többet. Mi az?"
$ws1.Range("A9").Value = "A földben születik, nagy tűzben égetik, karikára nyújtják, az ujjadra húzzák."
$ws1.Range("A10").Value = "Kis ember áll az erdőben, nagy kalap van a fejében. Ha kalapját megetted, őt magát is ismered."
$ws1.Range("A11").Value = "Kint is van, bent is van, mégis csak a házban van."
$ws1.Range("A12").Value = "Szekeremnek van kereke négy, de te azzal sehova se mégy. Én se ültem rajta soha még, messze van az ide, mint az ég."
$ws1.Range("A13").Value = "Lába van, de mégse jár, víz felett visz, nem madár. Nem rab, mégis láncot hord, s minden utat összetold."
$ws1.Range("A14").Value = "Kicsi is, görbe is. Tüzes még a feje is."
$ws1.Range("A15").Value = "Erdőn voltam, parton leszek, vízen járok, s tűzben veszek. Mi az?"
$ws1.Range("A16").Value = "Kis koromban nyersen esznek, ízét adom a levesnek. Barátom a petrezselyem, egy csomóba kötik velem."
$ws1.Range("A17").Value = "Zúgolódom, hangoskodom, pedig se szám, se pocakom, ami fogam alá kerül, megaprítom kegyetlenül, le azonban sosem nyelem, magam alá eresztgetem."
$ws1.Range("A18").Value = "Egy doboz, nincs se pántja, kulcsa vagy fedele, de aranyló kincset rejt a belseje."
$ws1.Range("A19").Value = "Mindent befal pofája, tátott szárnyast, szárnyatlant, fát, virágot, vasat csócsál, acélt ropogtat, fogával követ is kikoptat."
$ws1.Range("A20").Value = "Él lélegzettelen, halotti hidegen, sose szomjas, kortya örök, páncélt hord, mi sose zörög."
$ws1.Range("A21").Value = "Nem láthatni, nem tapinthatni, nem hallhatni, nem szagolhatni.Túl csillagokon, dombok tövében, kitölti az űrt egészen, sereghajtó, bár járt legelöl, életet végez, kacajt megöl."
$ws1.Range("A22").Value = "Mi az, minek a gyökere lappang, a fáknál magasabb, fel, fel, fel az égre tör, és mégis sose nő?"
$ws1.Range("A23").Value = "30 fehér ló egy piros dombon, abrakolnak, dobrokolnak majd rajtuk a béklyó."
$ws1.Range("A24").Value = "Hang nélkül kiált, szárnyatlan lebeg, fog nélkül kirág, szájatlan hebeg."

# --- Fill sheet2 (Valaszok) column A rows 1-24 with answers ---
$ws2.Range("A1").Value = "Válasz"
$ws2.Range("A2").Value = "szem"
$ws2.Range("A3").Value = "árnyék"
$ws2.Range("A4").Value = "nap"
$ws2.Range("A5").Value = "én"
$ws2.Range("A6").Value = "hagyma"
$ws2.Range("A7").Value = "szemüveg"
$ws2.Range("A8").Value = "nevem"
$ws2.Range("A9").Value = "gyűrű"
$ws2.Range("A10").Value = "gomba"
$ws2.Range("A11").Value = "ablak"
$ws2.Range("A12").Value = "göncölszekér"
$ws2.Range("A13").Value = "híd"
$ws2.Range("A14").Value = "pipa"
$ws2.Range("A15").Value = "csónak"
$ws2.Range("A16").Value = "répa"
$ws2.Range("A17").Value = "fűrész"
$ws2.Range("A18").Value = "tojás"
$ws2.Range("A19").Value = "idő"
$ws2.Range("A20").Value = "hal"
$ws2.Range("A21").Value = "sötétség"
$ws2.Range("A22").Value = "hegy"
$ws2.Range("A23").Value = "fogak"
$ws2.Range("A24").Value = "szél"

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 167
$ws2.Columns.Item(1).ColumnWidth = 39.35

# --- Sheet view / selection state ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws1.Range("A25").Select()

$ws2.Activate()
$ws2.Range("A13").Select()

$ws1.Activate()
